# Modify hlm dataloader workbook to load erosion data and soil data
# from the same Excel file: rename the original sheet to "erosion",
# add a new "soil" sheet with soil-sample data, and leave "soil" as
# the active/selected sheet (matching the commit's saved view state).

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet Sheet1 -> erosion -----------------------
$erosion = $wb.Worksheets.Item(1)
$erosion.Name = "erosion"

# Preserve the erosion sheet's own saved selection/scroll position
# (cursor left on C6 even though the sheet is no longer the active tab).
$erosion.Range("C6").Select()

# --- Add the new "soil" sheet right after "erosion" ---------------------
$soil = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $erosion)
$soil.Name = "soil"

# --- Header row -----------------------------------------------------------
$soil.Range("B1").Value = "PH值"
$soil.Range("C1").Value = "全氮含量(%)"
$soil.Range("D1").Value = "碳酸根离子"
$soil.Range("E1").Value = "硫酸根离子"
$soil.Range("F1").Value = "镁离子"
$soil.Range("G1").Value = "钾离子"
$soil.Range("H1").Value = "钠离子"
$soil.Range("I1").Value = "土壤电阻(Ω)"
$soil.Range("J1").Value = "站点"

# --- Data row -------------------------------------------------------------
$soil.Range("A2").Value = 0
$soil.Range("B2").Value = 6.75
$soil.Range("C2").Value = 0.109
$soil.Range("D2").Value = 0.0126
$soil.Range("E2").Value = 0.0118
$soil.Range("F2").Value = 0.0017
$soil.Range("G2").Value = 0.0004
$soil.Range("H2").Value = 0.0049
$soil.Range("I2").Value = 32.9
$soil.Range("J2").Value = "沈阳站"

# Leave the whole table selected on the soil sheet, and make "soil"
# the active (visible) tab when the workbook is opened -- matches the
# saved view (activeTab=1, tabSelected on the soil sheet).
$soil.Range("A1:J2").Select()
$soil.Activate()
